# Uncertainty_Tables_RD.xlsx - correct the numbering/condition for the
# "X_CAT = Average Distance" row on the Unified_table sheet.
#
# Row 43 (B43 = "L_E_12", C43 = "Average Distance") had a hard-coded index
# of 22 in column A, which duplicated the index already used by row 39 on
# the EnergyTransport sheet, breaking the running ID sequence on
# Unified_table. The rows below it (44-74) already compute their ID via
# formulas chained off A43 ("=A43+1", shared formulas, etc.), so fixing
# this single literal automatically cascades the correct sequential
# numbering (23, 24, 25, ... 54) through the remainder of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unified_table")
$ws.Activate()

$ws.Range("A43").Value = 23

# Restore the reviewer's last on-screen selection/scroll position for this
# sheet (it had moved down to inspect the corrected row and the rows that
# follow it).
$ws.Range("A44").Select()
